$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revise the count for 24/04/2020 (row 50, column B) from 1326 to 1408
$ws.Cells.Item(50, 2).Value = 1408

# Insert a new row for "25/04/2020" right after "25/03/2020" (row 51),
# pushing all subsequent rows down by one.
$ws.Cells.Item(52, 1).EntireRow.Insert()

# Populate the newly inserted row with the 25/04/2020 entry.
$ws.Cells.Item(52, 1).Value = "25/04/2020"
$ws.Cells.Item(52, 2).Value = 159
